$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(3,60259,9241,9894),
    @(4,36018,3284,4473),
    @(5,94724,8964,8487),
    @(6,1770,930,260),
    @(7,63763,10747,7884),
    @(8,7496,1462,1609),
    @(9,7549,1605,1046),
    @(10,3713,545,299),
    @(11,1403,399,3),
    @(12,0,0,0),
    @(13,1237,266,380),
    @(14,3607,1364,1349),
    @(15,6412,2261,1257),
    @(16,4529,2352,535),
    @(17,3616,1123,450),
    @(18,23732,3292,4483),
    @(19,1876,883,459),
    @(20,23821,3332,3720),
    @(21,461,575,129),
    @(22,23595,2657,4037),
    @(23,1574,835,259),
    @(24,26339,2686,4955),
    @(25,104631,8885,11887),
    @(26,8052,3060,1111),
    @(27,0,0,0),
    @(28,7028,1648,1725),
    @(29,3142,742,648),
    @(30,19640,3396,3866),
    @(31,620,128,379),
    @(32,3955,2269,409),
    @(33,18872,4474,3807),
    @(34,15225,4354,2774),
    @(35,6770,829,1347),
    @(36,76098,8585,7286),
    @(37,10775,3858,1490),
    @(38,36417,2610,3472),
    @(39,1591,1415,260),
    @(40,1782,668,713),
    @(41,3430,753,141),
    @(42,13802,693,358),
    @(43,364,265,58),
    @(44,1060,214,70),
    @(45,0,0,0),
    @(46,4428,1408,549),
    @(47,18728,4551,3571),
    @(48,45430,4453,6708),
    @(49,20982,4498,1923),
    @(50,16401,2040,3516),
    @(51,44303,4171,7019),
    @(52,6976,1309,1682),
    @(53,14930,2964,2463),
    @(54,2941,1818,1160),
    @(55,2928,1814,193),
    @(56,6977,1499,2663),
    @(57,15785,6243,3521),
    @(58,18656,1925,627),
    @(59,915772,138174,129216)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

[void]$ws.Range("B3").Select()
